# Apply the "Added new lora, added rgb leds, updated TSP65988 schematic,
# added bat to 3v3, added battery balancer" edit to the components sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header row: column D label changes from "doc" to "kogus" (quantity)
# ---------------------------------------------------------------
$ws.Range("D1").Value2 = "kogus"

# ---------------------------------------------------------------
# Row 7: LoRa module swapped for a new part (LoRa Ra-01SH)
# ---------------------------------------------------------------
$ws.Range("A7").Value2 = "LoRa Ra-01SH"
$ws.Range("B7").Value2 = "https://www.aliexpress.com/item/1005002702535557.html"
$ws.Range("C7").Value2 = 5.42

# ---------------------------------------------------------------
# Row 10: rename "USB PD" label to "USB-PD battery"
# ---------------------------------------------------------------
$ws.Range("A10").Value2 = "USB-PD battery"

# ---------------------------------------------------------------
# Row 11 (USB-C): price cell is cleared (no longer priced individually)
# quantity of 2 is introduced
# ---------------------------------------------------------------
$ws.Range("C11").ClearContents()

# ---------------------------------------------------------------
# Row 12 / 13: Speaker amplifier / Piezo speaker keep their labels,
# but lose any link/price they might have had (they stay bare).
# (no value changes needed besides quantity/E formulas added later)
# ---------------------------------------------------------------

# ---------------------------------------------------------------
# Row 15: "Levelshifter/isolator" replaced by "RGB LED"
# ---------------------------------------------------------------
$ws.Range("A15").Value2 = "RGB LED"
$ws.Range("B15").Value2 = "https://www.aliexpress.com/item/32453497583.html"

# ---------------------------------------------------------------
# Row 16: "IO connectors" entry removed entirely
# ---------------------------------------------------------------
$ws.Range("A16").ClearContents()

# ---------------------------------------------------------------
# Row 18 (new): USB - UART converter
# ---------------------------------------------------------------
$ws.Range("A18").Value2 = "USB - UART "
$ws.Range("B18").Value2 = "https://www.mouser.ee/ProductDetail/Silicon-Labs/CP2102N-A02-GQFN20?qs=u16ybLDytRaG8WdlP0fT2g%3D%3D"
$ws.Range("C18").Value2 = 2.16

# ---------------------------------------------------------------
# Row 19: previously only had the bare TPS65988 link in column A;
# now gets a proper label in A, the link moves to B, and a price.
# ---------------------------------------------------------------
$ws.Range("A19").Value2 = "USB-PD negotiation"
$ws.Range("B19").Value2 = "https://www.mouser.ee/ProductDetail/Texas-Instruments/TPS65988DJRSHR?qs=T3oQrply3y8cml9f5FMm3A%3D%3D"
$ws.Range("C19").Value2 = 6.64

# ---------------------------------------------------------------
# Row 21 (new): Digi pote
# ---------------------------------------------------------------
$ws.Range("A21").Value2 = "Digi pote"
$ws.Range("B21").Value2 = "https://www.mouser.ee/ProductDetail/Microchip-Technology-Atmel/MCP40D18T-503E-LT?qs=dQMF8gqycOVbxabKF9CRwQ%3D%3D"
$ws.Range("C21").Value2 = 0.69
$ws.Range("D21").Value2 = 2

# ---------------------------------------------------------------
# Row 22 (new): Bat to USB-c
# ---------------------------------------------------------------
$ws.Range("A22").Value2 = "Bat to USB-c"
$ws.Range("B22").Value2 = "https://www.mouser.ee/ProductDetail/Monolithic-Power-Systems-MPS/MP2229GQ-P?qs=ZNK0BnemlqHKDD1LC56W1w%3D%3D"
$ws.Range("C22").Value2 = 2.84
$ws.Range("D22").Value2 = 2

# ---------------------------------------------------------------
# Row 23 (new): Bat balancer
# ---------------------------------------------------------------
$ws.Range("A23").Value2 = "Bat balancer"
$ws.Range("B23").Value2 = "https://www.mouser.ee/ProductDetail/Texas-Instruments/BQ29209DRBR?qs=hEBn5lgDlCoqdeLwAzko8w%3D%3D"
$ws.Range("C23").Value2 = 0.96
$ws.Range("D23").Value2 = 1

# ---------------------------------------------------------------
# New "kogus" (quantity) column D - set to 1 for every existing row
# that carries a price, 0 for the battery row (not ordered), and 2
# for USB-C / already-set rows above.
# ---------------------------------------------------------------
$ws.Range("D2").Value2 = 1
$ws.Range("D3").Value2 = 1
$ws.Range("D4").Value2 = 1
$ws.Range("D5").Value2 = 1
$ws.Range("D6").Value2 = 1
$ws.Range("D7").Value2 = 1
$ws.Range("D8").Value2 = 1
$ws.Range("D9").Value2 = 0
$ws.Range("D10").Value2 = 1
$ws.Range("D11").Value2 = 2
$ws.Range("D14").Value2 = 1
$ws.Range("D18").Value2 = 1
$ws.Range("D19").Value2 = 1

# ---------------------------------------------------------------
# New "kokku" (total) column E = price * quantity, for every row
# from 2 through 28. E2 is its own formula; E3:E28 is filled as one
# block so Excel records it as a single shared formula (matching
# how the workbook was actually edited).
# ---------------------------------------------------------------
$ws.Range("E2").Formula = "=C2*D2"
$ws.Range("E3:E28").Formula = "=C3*D3"

# ---------------------------------------------------------------
# Grand total formula now sums the new E column instead of C.
# ---------------------------------------------------------------
$ws.Range("G1").Formula = "=SUM(E:E)"

# ---------------------------------------------------------------
# Selection / view state, best effort.
# ---------------------------------------------------------------
$ws.Range("B23").Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 16
    $win.ScrollColumn = 1
}
